$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracker rows for 2025-09-14 (serial date 45914), appended after the
# existing data (last row was 36, date 45913).
$newRows = @(
    @("G2", "Workout", 45914, 0.9327180547071353, 0, -0.01),
    @("G3", "Eat Healthy", 45914, 0.9327180547071353, 0, -0.01),
    @("G4", "Read Book", 45914, 0.9327180547071353, 0, -0.01),
    @("G5", "Investment Plan", 45914, 0.9327180547071353, 0, -0.01),
    @("G6", "Spend 10 Hours without phone", 45914, 0.9327180547071353, 0, -0.01)
)

$startRow = 37
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]

    $dateCell = $ws.Cells.Item($r, 3)
    $dateCell.Value = $data[2]
    $dateCell.NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
